$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# --- Move DepoYeri_2..4 rows (16-18) down to rows 22-24, row 15 becomes the active entry ---

# Capture current text/values before we start overwriting anything.
$a16 = $ws.Range("A16").Value2
$b16 = $ws.Range("B16").Value2
$a17 = $ws.Range("A17").Value2
$b17 = $ws.Range("B17").Value2
$a18 = $ws.Range("A18").Value2
$b18 = $ws.Range("B18").Value2

# Row 15 - value changes from 3001 to 3002 (now the active DepoYeri)
$ws.Range("B15").Value = 3002

# Clear old rows 16-18 (content moves to 22-24)
$ws.Range("A16:B18").ClearContents()

# Write moved content into rows 22-24
$ws.Range("A22").Value = $a16
$ws.Range("B22").Value = $b16
$ws.Range("A23").Value = $a17
$ws.Range("B23").Value = $b17
$ws.Range("A24").Value = $a18
$ws.Range("B24").Value = $b18

# Re-apply the number-style that B16/B17/B18 used to B22/B23/B24 (matches cellXfs idx 8 pattern: font+left align)
$ws.Range("B22:B24").HorizontalAlignment = -4131

# --- Selection update ---
$ws.Range("B16").Select()
